$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.StartsWith("System, ")) {
        $parts = $val.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        if ($trimmed[0] -eq "System") {
            $rest = @()
            for ($i = 1; $i -lt $trimmed.Count; $i++) {
                $rest += $trimmed[$i]
            }
            $rest += "System"
            $newVal = [string]::Join(", ", $rest)
            $cell.Value2 = $newVal
        }
    }
}
